# Applies the text edits described by the commit diff to the active
# document using Word's Find/Replace (wdFindContinue=1, wdReplaceOne=2).

$d = $word.ActiveDocument

function Replace-Text($findText, $replaceText) {
    $d.Content.Find.Execute($findText, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $replaceText, 2)
}

# 1. "мог бы посчитать" -> "мог посчитать"
Replace-Text "мог бы посчитать, что находится перед домом лорда в Лондоне." `
             "мог посчитать, что находится перед домом лорда в Лондоне."

# 2. "состоятельным коллегам на Земле" -> "состоятельным товарищам на Земле"
Replace-Text "состоятельным коллегам на Земле" "состоятельным товарищам на Земле"

# 3. ", способствовать" -> " способствовать" (comma -> space before the word)
Replace-Text "удобно, способствовать" "удобно способствовать"

# 4. "Машины из стали и латуни" -> "Механизмы из стали и латуни"
Replace-Text " Машины из стали и латуни, винт" " Механизмы из стали и латуни, винт"

# 5. ", от бронированных пластин до привода" -> ", бронированные пластины"
Replace-Text "шестерни, от бронированных пластин до привода" "шестерни, бронированные пластины"

# 6. "завод будет выгоден." -> "завод будет не лишним."
Replace-Text "металлообрабатывающий завод будет выгоден." `
             "металлообрабатывающий завод будет не лишним."

# 7. "металлический кабель высоко" -> "металлический трос высоко"
Replace-Text "металлический кабель высоко" "металлический трос высоко"

# 8a. Drop the English parentheticals after the station names.
Replace-Text "Горький Пролом (Sourbreak), Бурелом (Deadfall) и Вид на реку (Riverview)." `
             "Горький Пролом, Бурелом и Вид на реку."

# 8b. "может стать смертным приговором." -> "может стоить жизни."
Replace-Text "может стать смертным приговором." "может стоить жизни."

# 8c. add a comma after "кораблю"
Replace-Text "открыла огонь по кораблю пролетающему над" `
             "открыла огонь по кораблю, пролетающему над"

# 8d. reorder "слишком опасно совершать полеты вне города" -> "совершать полеты вне города слишком опасно"
Replace-Text "Генерал-Губернатор решил, что слишком опасно совершать полеты вне города." `
             "Генерал-Губернатор решил, что совершать полеты вне города слишком опасно."
